$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal text value into a cell without changing its
# effective style (avoids Excel auto-converting numeric-looking strings
# such as "1.00" or "0.998" into actual numbers).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.525.81"
$ws.Range("E2").Value = "  +0.48%  "
Set-TextValue $ws.Range("D3") "3.184.94"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "601.51"
$ws.Range("E5").Value = "  -0.69%  "
Set-TextValue $ws.Range("D6") "155.09"
$ws.Range("E6").Value = "  -0.10%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue $ws.Range("D8") "3.183.27"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("E10").Value = "  -2.14%  "
Set-TextValue $ws.Range("D11") "5.63"
$ws.Range("E11").Value = "  -9.18%  "
$ws.Range("E12").Value = "  +0.61%  "
Set-TextValue $ws.Range("D13") "0.0000268"
$ws.Range("E13").Value = "  -1.37%  "
Set-TextValue $ws.Range("D14") "38.76"
$ws.Range("E14").Value = "  -0.54%  "
Set-TextValue $ws.Range("D15") "3.707.34"
$ws.Range("E15").Value = "  -1.00%  "
Set-TextValue $ws.Range("D16") "66.514.09"
$ws.Range("E16").Value = "  +0.26%  "
Set-TextValue $ws.Range("D17") "7.44"
$ws.Range("E17").Value = "  -0.90%  "
Set-TextValue $ws.Range("D18") "3.182.64"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("E19").Value = "  +0.42%  "
Set-TextValue $ws.Range("D20") "513.91"
$ws.Range("E20").Value = "  -0.22%  "
Set-TextValue $ws.Range("D21") "15.47"
$ws.Range("E21").Value = "  -2.39%  "
Set-TextValue $ws.Range("D22") "0.734"
$ws.Range("E22").Value = "  -0.78%  "
Set-TextValue $ws.Range("D23") "8.15"
$ws.Range("E23").Value = "  +1.60%  "
Set-TextValue $ws.Range("D24") "14.92"
$ws.Range("E24").Value = "  -2.60%  "
Set-TextValue $ws.Range("D25") "84.82"
$ws.Range("E25").Value = "  -1.04%  "
Set-TextValue $ws.Range("D26") "0.998"
$ws.Range("E26").Value = "  -0.05%  "
Set-TextValue $ws.Range("D27") "9.25"
$ws.Range("E27").Value = "  -0.38%  "
Set-TextValue $ws.Range("D28") "3.00"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("E29").Value = "  +6.75%  "
Set-TextValue $ws.Range("D30") "3.13"
$ws.Range("E30").Value = "  +7.99%  "
$ws.Range("E31").Value = "  +1.48%  "
Set-TextValue $ws.Range("D32") "28.12"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  -1.41%  "
Set-TextValue $ws.Range("D35") "6.55"
$ws.Range("E35").Value = "  -2.21%  "
Set-TextValue $ws.Range("D36") "515.53"
$ws.Range("E36").Value = "  +4.62%  "
Set-TextValue $ws.Range("D37") "54.83"
Set-TextValue $ws.Range("D38") "0.0891"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  +6.27%  "
Set-TextValue $ws.Range("D41") "8.88"
$ws.Range("E41").Value = "  -0.20%  "
Set-TextValue $ws.Range("D42") "0.0₃0683"
$ws.Range("E42").Value = "  +5.00%  "
Set-TextValue $ws.Range("D43") "0.300"
$ws.Range("E43").Value = "  +1.18%  "
Set-TextValue $ws.Range("D44") "2.80"
$ws.Range("E44").Value = "  -8.13%  "
Set-TextValue $ws.Range("D45") "2.46"
$ws.Range("E45").Value = "  -2.05%  "
Set-TextValue $ws.Range("D46") "2.851.14"
$ws.Range("E46").Value = "  -5.58%  "
Set-TextValue $ws.Range("D47") "28.20"
$ws.Range("E47").Value = "  -3.73%  "
Set-TextValue $ws.Range("D48") "2.39"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D49") "0.999"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D50") "0.117"
$ws.Range("E50").Value = "  +0.77%  "
Set-TextValue $ws.Range("D51") "2.58"
$ws.Range("E51").Value = "  +5.64%  "
